# Update launch parameters with Space Shuttle options
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Launch Parameters")

# Insert a new row above row 4, pushing the Saturn V..Atlas V 551 rows down
# (and dropping the old Space Shuttle row that used to sit at row 9, since
# the table keeps a fixed 9-row extent).
$ws.Rows.Item(4).Insert()

# New Space Shuttle launch parameters go into the freshly inserted row 4.
# (Row-insert already carries the neighbouring rows' per-column formatting,
# matching the workbook's existing s="2"/s="3" styling, so no extra
# formatting step is needed here.)
$ws.Cells.Item(4, 1).Value = "Space Shuttle"
$ws.Cells.Item(4, 2).Value = 375
$ws.Cells.Item(4, 3).Value = 0.25
$ws.Cells.Item(4, 4).Value = 50
$ws.Cells.Item(4, 5).Value = 165
$ws.Cells.Item(4, 6).Value = 8
$ws.Cells.Item(4, 7).Value = 40
$ws.Cells.Item(4, 8).Value = 180
$ws.Cells.Item(4, 9).Value = 180

# Drop the now-duplicated last row (old row 9 / Space Shuttle entry shifted
# to row 10 by the insert) so the table stays at 9 data rows.
$ws.Rows.Item(10).Delete()

# Match the new selection recorded in the sheet view.
$ws.Range("G4").Select()
